$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.735883
$ws.Range("H2").Value = 11.207649
$ws.Range("I2").Value = 0.1699536238627456
$ws.Range("J2").Value = 0.1699536238627456
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4794123333333333
$ws.Range("N2").Value = 1.438237
$ws.Range("O2").Value = 0.2944679249717814
$ws.Range("P2").Value = 0.2944679249717814
$ws.Range("Q2").Value = 1.791028386090333
$ws.Range("R2").Value = 16.119255474813
$ws.Range("S2").Value = 0.05004589096029732
$ws.Range("T2").Value = 0.05004589096029732
$ws.Range("G3").Value = 3.735883
$ws.Range("H3").Value = 11.207649
$ws.Range("I3").Value = 0.1699536238627456
$ws.Range("J3").Value = 0.1699536238627456
$ws.Range("O3").Value = 0.4529251427412002
$ws.Range("P3").Value = 0.4529251427412002
$ws.Range("Q3").Value = 2.754805255958666
$ws.Range("R3").Value = 24.793247303628
$ws.Range("S3").Value = 0.0769762693474183
$ws.Range("T3").Value = 0.07697626934741827
$ws.Range("G4").Value = 3.735883
$ws.Range("H4").Value = 11.207649
$ws.Range("I4").Value = 0.1699536238627456
$ws.Range("J4").Value = 0.1699536238627456
$ws.Range("M4").Value = 0.41126
$ws.Range("N4").Value = 1.23378
$ws.Range("O4").Value = 0.2526069322870184
$ws.Range("P4").Value = 0.2526069322870184
$ws.Range("Q4").Value = 1.53641924258
$ws.Range("R4").Value = 13.82777318322
$ws.Range("S4").Value = 0.04293146355502996
$ws.Range("T4").Value = 0.04293146355502997
$ws.Range("I5").Value = 0.620557022856673
$ws.Range("J5").Value = 0.6205570228566729
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4794123333333333
$ws.Range("N5").Value = 1.438237
$ws.Range("O5").Value = 0.2944679249717814
$ws.Range("P5").Value = 0.2944679249717814
$ws.Range("Q5").Value = 6.539638389950445
$ws.Range("R5").Value = 58.856745509554
$ws.Range("S5").Value = 0.1827341388472708
$ws.Range("T5").Value = 0.1827341388472708
$ws.Range("I6").Value = 0.620557022856673
$ws.Range("J6").Value = 0.6205570228566729
$ws.Range("O6").Value = 0.4529251427412002
$ws.Range("P6").Value = 0.4529251427412002
$ws.Range("Q6").Value = 10.05870724809156
$ws.Range("R6").Value = 90.528365232824
$ws.Range("S6").Value = 0.2810658781564129
$ws.Range("T6").Value = 0.2810658781564128
$ws.Range("I7").Value = 0.620557022856673
$ws.Range("J7").Value = 0.6205570228566729
$ws.Range("M7").Value = 0.41126
$ws.Range("N7").Value = 1.23378
$ws.Range("O7").Value = 0.2526069322870184
$ws.Range("P7").Value = 0.2526069322870184
$ws.Range("Q7").Value = 5.609976000306667
$ws.Range("R7").Value = 50.48978400276
$ws.Range("S7").Value = 0.1567570058529893
$ws.Range("T7").Value = 0.1567570058529893
$ws.Range("G8").Value = 4.604948666666666
$ws.Range("H8").Value = 13.814846
$ws.Range("I8").Value = 0.2094893532805814
$ws.Range("J8").Value = 0.2094893532805814
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.4794123333333333
$ws.Range("N8").Value = 1.438237
$ws.Range("O8").Value = 0.2944679249717814
$ws.Range("P8").Value = 0.2944679249717814
$ws.Range("Q8").Value = 2.207669185166889
$ws.Range("R8").Value = 19.869022666502
$ws.Range("S8").Value = 0.06168789516421325
$ws.Range("T8").Value = 0.06168789516421325
$ws.Range("G9").Value = 4.604948666666666
$ws.Range("H9").Value = 13.814846
$ws.Range("I9").Value = 0.2094893532805814
$ws.Range("J9").Value = 0.2094893532805814
$ws.Range("O9").Value = 0.4529251427412002
$ws.Range("P9").Value = 0.4529251427412002
$ws.Range("Q9").Value = 3.395646167279111
$ws.Range("R9").Value = 30.560815505512
$ws.Range("S9").Value = 0.09488299523736907
$ws.Range("T9").Value = 0.09488299523736904
$ws.Range("G10").Value = 4.604948666666666
$ws.Range("H10").Value = 13.814846
$ws.Range("I10").Value = 0.2094893532805814
$ws.Range("J10").Value = 0.2094893532805814
$ws.Range("M10").Value = 0.41126
$ws.Range("N10").Value = 1.23378
$ws.Range("O10").Value = 0.2526069322870184
$ws.Range("P10").Value = 0.2526069322870184
$ws.Range("Q10").Value = 1.893831188653333
$ws.Range("R10").Value = 17.04448069788
$ws.Range("S10").Value = 0.0529184628789991
$ws.Range("T10").Value = 0.05291846287899911
